# Applies the "Trade #54 closed" update to the live trading results workbook:
#  - Updates rolled-up summary figures on the "Summary" sheet
#  - Updates the MarketMaking strategy row on the "Strategy Status" sheet
#  - Appends the new closed trade as row 55 on both the "All Trades" and
#    "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.33   # Current Capital
$summary.Range("B4").Value = 0.33      # Total P&L $
$summary.Range("B5").Value = 0.12      # Total P&L %
$summary.Range("B6").Value = 54        # Total Trades
$summary.Range("B8").Value = 30        # Losing Trades
$summary.Range("B9").Value = 27.78     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.33     # Capital
$status.Range("D4").Value = 54         # Trades
$status.Range("E4").Value = 0.33       # P&L $
$status.Range("F4").Value = 0.33       # P&L %
$status.Range("G4").Value = 27.78      # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade row (row 55) to a trades sheet
# ---------------------------------------------------------------------------
function Add-Trade54Row($ws) {
    $ws.Range("A55").Value = 54

    # B55 looks like a date ("2026-02-17") and would otherwise be silently
    # auto-converted into a date serial by the COM layer, whereas the source
    # workbook keeps it as plain text. Force text formatting just for this
    # assignment, then restore the default "Normal" style so the cell is
    # left without any explicit style override (matching the rest of the
    # sheet).
    $ws.Range("B55").NumberFormat = "@"
    $ws.Range("B55").Value = "2026-02-17"
    $ws.Range("B55").Style = "Normal"

    $ws.Range("C55").Value = "15:43:02"
    $ws.Range("D55").Value = "MarketMaking"
    $ws.Range("E55").Value = "DOWN"
    $ws.Range("F55").Value = 0.77
    $ws.Range("G55").Value = 0.68
    $ws.Range("H55").Value = "CLOSED"
    $ws.Range("I55").Value = -11.6883
    $ws.Range("J55").Value = -0.09
    $ws.Range("K55").Value = 100.33
    $ws.Range("L55").Value = 0
    $ws.Range("M55").Value = 0
    $ws.Range("N55").Value = 0.6
    $ws.Range("O55").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P55").Value = "early_exit"
    $ws.Range("Q55").Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade54Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade54Row $marketMaking

Write-Output "Applied trade #54 update"
